$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) Insert a brand-new red-highlighted bullet BEFORE the very first paragraph:
#    "Report05 + 06 (Installations) goed krijgen."
# ---------------------------------------------------------------------------
$firstPara = $d.Paragraphs(1)
$firstPara.Range.InsertParagraphBefore()

$reportFrag = "<w:p $wns>" +
    "<w:pPr>" +
        "<w:pStyle w:val=`"Lijstalinea`"/>" +
        "<w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr>" +
        "<w:rPr><w:highlight w:val=`"red`"/></w:rPr>" +
    "</w:pPr>" +
    "<w:r><w:rPr><w:highlight w:val=`"red`"/></w:rPr><w:t>Report05 + 06 (</w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:rPr><w:highlight w:val=`"red`"/></w:rPr><w:t>Installations</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:rPr><w:highlight w:val=`"red`"/></w:rPr><w:t>) goed krijgen.</w:t></w:r>" +
    "</w:p>"
$d.Paragraphs(1).Range.InsertXML($reportFrag)

# ---------------------------------------------------------------------------
# 2) "Fix release nummers in namen (REGEX template toevoegen)" paragraph:
#    highlight changes from yellow to green (paragraph mark + all 3 runs)
# ---------------------------------------------------------------------------
$fixFrag = "<w:p $wns>" +
    "<w:pPr>" +
        "<w:pStyle w:val=`"Lijstalinea`"/>" +
        "<w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr>" +
        "<w:rPr><w:highlight w:val=`"green`"/></w:rPr>" +
    "</w:pPr>" +
    "<w:r><w:rPr><w:highlight w:val=`"green`"/></w:rPr><w:t>Fix release nummers in namen (REGEX template toevoegen)</w:t></w:r>" +
    "<w:r><w:rPr><w:highlight w:val=`"green`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r>" +
    "<w:r><w:rPr><w:highlight w:val=`"green`"/></w:rPr><w:sym w:font=`"Wingdings`" w:char=`"F0E7`"/></w:r>" +
    "</w:p>"

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.StartsWith("Fix release nummers in namen")) {
        $p.Range.InsertXML($fixFrag)
        break
    }
}

# ---------------------------------------------------------------------------
# 3) "Locatie in lijsten onthouden" paragraph: add yellow highlight
#    (paragraph mark + run) where there was none before
# ---------------------------------------------------------------------------
$locatieFrag = "<w:p $wns>" +
    "<w:pPr>" +
        "<w:pStyle w:val=`"Lijstalinea`"/>" +
        "<w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr>" +
        "<w:rPr><w:highlight w:val=`"yellow`"/></w:rPr>" +
    "</w:pPr>" +
    "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>Locatie in lijsten onthouden</w:t></w:r>" +
    "</w:p>"

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.StartsWith("Locatie in lijsten onthouden")) {
        $p.Range.InsertXML($locatieFrag)
        break
    }
}

# ---------------------------------------------------------------------------
# 4) Append two new yellow-highlighted bullets at the very end of the document
#    (after "... sensor herschrijven", before the sectPr)
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$idsFrag = "<w:p $wns>" +
    "<w:pPr>" +
        "<w:pStyle w:val=`"Lijstalinea`"/>" +
        "<w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr>" +
        "<w:rPr><w:highlight w:val=`"yellow`"/></w:rPr>" +
    "</w:pPr>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>ID’s</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t xml:space=`"preserve`"> tonen in diverse schermen</w:t></w:r>" +
    "</w:p>"
$d.Paragraphs($d.Paragraphs.Count).Range.InsertXML($idsFrag)

$lastPara2 = $d.Paragraphs($d.Paragraphs.Count)
$lastPara2.Range.InsertParagraphAfter()
$overviewFrag = "<w:p $wns>" +
    "<w:pPr>" +
        "<w:pStyle w:val=`"Lijstalinea`"/>" +
        "<w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr>" +
        "<w:rPr><w:highlight w:val=`"yellow`"/></w:rPr>" +
    "</w:pPr>" +
    "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t xml:space=`"preserve`">OVERVIEW rapport </w:t></w:r>" +
    "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:sym w:font=`"Wingdings`" w:char=`"F0E8`"/></w:r>" +
    "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>Vendor</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t xml:space=`"preserve`"> – Component - Installaties</w:t></w:r>" +
    "</w:p>"
$d.Paragraphs($d.Paragraphs.Count).Range.InsertXML($overviewFrag)

Write-Host "edit complete"
